$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M4").Value = "Сидоров"
$ws.Range("M5").Value = "групи 117-а"
$ws.Range("A8").Value = "без порушень встановлених термінів і позитивної підсумкової семестрової модульної рейтингової оцінки,`nпросимо Вашого дозволу на звільнення нас від складання семестрового екзамену з даної дисципліни`nпровідний викладач Вечерковська А. С. і зарахування відповідної`nпідсумкової семестрової рейтингової оцінки."
$ws.Range("A9").Value = "У зв’язку з отриманням у 1-му семестрі 2019/2020 навчального року позитивних підсумкових`nмодульних рейтингових оцінок з усіх 5 модулів дисципліни`nОснови ООП"
$ws.Range("B10").Value = "Прізвище та ініціали`nстудента"
$ws.Range("M10").Value = "Підпис`nстудента"
$ws.Range("C11").Value = "Підсумкова`nмодульна (бали)"
$ws.Range("F11").Value = "Підсумкова`nсеместрова`nмодульна`nрейтингова оцінка"
$ws.Range("J11").Value = "Підсумкова`nсеместрова`nрейтингова оцінка"
$ws.Range("C12").Value = "Мод.`n№1"
$ws.Range("D12").Value = "Мод.`n№2"
$ws.Range("E12").Value = "Мод.`n№3"
$ws.Range("G12").Value = "Націон.`nшкала"
$ws.Range("I12").Value = "Націон.`nшкала"
$ws.Range("K12").Value = "Націон.`nшкала"
$ws.Range("L12").Value = "Шкала`nECTS"
$ws.Range("A13").Value = "1"
$ws.Range("B13").Value = "Галацюк Т. П."
$ws.Range("C13").Value = "test"
$ws.Range("D13").Value = "68"
$ws.Range("E13").Value = "57"
$ws.Range("F13").Value = "0b"
$ws.Range("G13").Value = "tn890"
$ws.Range("H13").Value = "n87"
$ws.Range("I13").Value = "5tn8b0"
$ws.Range("J13").Value = "kjh"
$ws.Range("K13").Value = "o8"
$ws.Range("L13").Value = "tuyf"
$ws.Range("M13").Value = "978rf"
$ws.Range("A14").Value = "2"
$ws.Range("B14").Value = "Батрак О. П."
$ws.Range("C14").Value = "87"
$ws.Range("D14").Value = "8f"
$ws.Range("E14").Value = "8g"
$ws.Range("F14").Value = "sdfg"
$ws.Range("G14").Value = "8g"
$ws.Range("H14").Value = "8"
$ws.Range("I14").Value = "tn"
$ws.Range("J14").Value = "89"
$ws.Range("K14").Value = "n"
$ws.Range("L14").Value = "r7"
$ws.Range("M14").Value = "r"
